# Apply "New crime data collected" update to the CompStat weekly report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text updates (Volume/Number line and Report-covering-week line).
#    These are rich-text shared strings made up of several runs; we update
#    the whole caption text (formatting for the cell as a whole is retained).
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/22/2025  Through  12/28/2025"

# ---------------------------------------------------------------------------
# 2) Plain numeric value refreshes across the crime-complaints table
#    (rows 14-31), columns C..N.
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -88.888888888888

$ws.Range("C15").Value = 1
$ws.Range("G15").Value = 7
$ws.Range("H15").Value = -14.285714285714
$ws.Range("I15").Value = 84
$ws.Range("K15").Value = 5
$ws.Range("L15").Value = 35.483870967741
$ws.Range("M15").Value = 37.704918032786
$ws.Range("N15").Value = -30

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 44
$ws.Range("G16").Value = 37
$ws.Range("H16").Value = 18.918918918918
$ws.Range("I16").Value = 528
$ws.Range("J16").Value = 634
$ws.Range("K16").Value = -16.719242902208
$ws.Range("L16").Value = -22.919708029197
$ws.Range("M16").Value = -30.158730158730
$ws.Range("N16").Value = -83.060635226179

$ws.Range("C17").Value = 19
$ws.Range("E17").Value = 18.75
$ws.Range("F17").Value = 67
$ws.Range("G17").Value = 66
$ws.Range("H17").Value = 1.515151515151
$ws.Range("I17").Value = 1054
$ws.Range("J17").Value = 1077
$ws.Range("K17").Value = -2.135561745589
$ws.Range("L17").Value = -2.135561745589
$ws.Range("M17").Value = 42.240215924426
$ws.Range("N17").Value = -28.250510551395

$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 22
$ws.Range("H18").Value = 46.666666666666
$ws.Range("I18").Value = 245
$ws.Range("J18").Value = 281
$ws.Range("K18").Value = -12.811387900355
$ws.Range("L18").Value = -30.790960451977
$ws.Range("M18").Value = -48.093220338983
$ws.Range("N18").Value = -86.655773420479

$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -30.769230769230
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = 23.809523809523
$ws.Range("I19").Value = 905
$ws.Range("J19").Value = 908
$ws.Range("K19").Value = -0.330396475770
$ws.Range("L19").Value = -16.589861751152
$ws.Range("M19").Value = 27.464788732394
$ws.Range("N19").Value = -5.433646812957

$ws.Range("C20").Value = 17
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 240
$ws.Range("F20").Value = 44
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 46.666666666666
$ws.Range("I20").Value = 481
$ws.Range("J20").Value = 594
$ws.Range("K20").Value = -19.023569023569
$ws.Range("L20").Value = -12.386156648451
$ws.Range("M20").Value = 49.844236760124
$ws.Range("N20").Value = -81.787201817493

$ws.Range("C21").Value = 57
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = 39.024390243902
$ws.Range("F21").Value = 236
$ws.Range("G21").Value = 198
$ws.Range("H21").Value = 19.191919191919
$ws.Range("I21").Value = 3311
$ws.Range("J21").Value = 3586
$ws.Range("K21").Value = -7.668711656441
$ws.Range("L21").Value = -13.663624511082
$ws.Range("M21").Value = 7.082794307891
$ws.Range("N21").Value = -67.747905708162

$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("L22").Value = -11.538461538461
$ws.Range("M22").Value = -38.666666666666

$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 15
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 282
$ws.Range("J23").Value = 280
$ws.Range("K23").Value = 0.714285714285
$ws.Range("L23").Value = -13.761467889908
$ws.Range("M23").Value = 24.229074889867

$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 59
$ws.Range("E24").Value = -44.067796610169
$ws.Range("F24").Value = 191
$ws.Range("G24").Value = 213
$ws.Range("H24").Value = -10.328638497652
$ws.Range("I24").Value = 2591
$ws.Range("J24").Value = 2433
$ws.Range("K24").Value = 6.494040279490
$ws.Range("L24").Value = 12.996075010902
$ws.Range("M24").Value = 72.273936170212

$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -42.105263157894
$ws.Range("F25").Value = 60
$ws.Range("G25").Value = 73
$ws.Range("H25").Value = -17.808219178082
$ws.Range("I25").Value = 909
$ws.Range("J25").Value = 939
$ws.Range("K25").Value = -3.194888178913
$ws.Range("L25").Value = 1

$ws.Range("C26").Value = 39
$ws.Range("D26").Value = 32
$ws.Range("E26").Value = 21.875
$ws.Range("F26").Value = 109
$ws.Range("G26").Value = 108
$ws.Range("H26").Value = 0.925925925925
$ws.Range("I26").Value = 1467
$ws.Range("J26").Value = 1541
$ws.Range("K26").Value = -4.802076573653
$ws.Range("L26").Value = 16.336241078509
$ws.Range("M26").Value = -19.660460021905

$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -12.5
$ws.Range("I27").Value = 98
$ws.Range("K27").Value = -4.854368932038
$ws.Range("L27").Value = 1.030927835051

$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 14
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 133.333333333333
$ws.Range("I28").Value = 140
$ws.Range("J28").Value = 135
$ws.Range("K28").Value = 3.703703703703
$ws.Range("L28").Value = 21.739130434782

$ws.Range("L29").Value = -25
$ws.Range("N29").Value = -84.615384615384

$ws.Range("L30").Value = -35.185185185185
$ws.Range("N30").Value = -87.632508833922

$ws.Range("L31").Value = -62.5

# ---------------------------------------------------------------------------
# 3) Cells that flip between a numeric value and the "not enough data" text
#    placeholders ("0" = shared text string index 20, "***.*" = index 21),
#    mirroring the styling already used on row 14 for the same situation.
#    We set the textual value first (prefixed with an apostrophe so it is
#    stored as text rather than being auto-converted back to a number) and
#    then copy over the number-formatting/style from a matching reference
#    cell so the style index lines up with the rest of the "N/A" cells.
# ---------------------------------------------------------------------------

# Row 22: C22 & D22 -> "0" (text), E22 -> "***.*" (text), styled like row 14.
$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

# Row 27: D27 -> "0" (text), E27 -> "***.*" (text), styled like row 14.
$ws.Range("D27").Value = "'0"
$ws.Range("D14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# Row 28: D28 & E28 switch back from the text placeholders to real numbers,
# restyled like their numeric neighbours (C28 = number style, H28 = percent style).
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 4

$ws.Range("H28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -75
